# STM32F10x Pin Assignment workbook edit
# Commit: "Swapped PB4 with PB5 for I2C alert functionality."
#
# Row 42 holds pin PB4 (I2C1_SMBA alt-function note) and row 43 holds pin PB5.
# The "Connected to" / "Description" / "Comments" values (columns H, I, J)
# that used to describe PB4's SmartFusion/Programmer-header/JTAG routing
# actually belong to PB5, and the battery-gauge (LTC2942-1 / I2C #AL/CC alert)
# routing that used to be listed under PB5 actually belongs to PB4.
# This swaps those three columns between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$null = $ws.Activate()

# Capture the current ("before") values for row 42 (PB4) and row 43 (PB5)
$h42 = $ws.Range("H42").Value2
$i42 = $ws.Range("I42").Value2
$j42 = $ws.Range("J42").Value2
$h43 = $ws.Range("H43").Value2
$i43 = $ws.Range("I43").Value2
$j43 = $ws.Range("J43").Value2

# Swap them: row 42 gets what used to be in row 43, and vice versa
$ws.Range("H42").Value2 = $h43
$ws.Range("I42").Value2 = $i43
$ws.Range("J42").Value2 = $j43
$ws.Range("H43").Value2 = $h42
$ws.Range("I43").Value2 = $i42
$ws.Range("J43").Value2 = $j42

# Row 43 now contains the long wrapped "SmartFusion connector | Programmer
# header" text (formerly on row 42), so it needs the taller row height that
# row 42 already had; row 42 keeps its own height (driven by other wrapped
# cells in that row) so it is left untouched.
$ws.Rows.Item(43).RowHeight = 30

# Restore a plausible selection/view state similar to the saved workbook.
$null = $ws.Range("B36").Select()
